$d = $word.ActiveDocument
$sec = $d.Sections.First

# The document has two distinct headers (primary + first-page) and two
# distinct footers (primary + first-page), each containing exactly one
# inline picture (logo). Per the commit, the logos' internal "name"
# identifiers are swapped:
#   Pearson logo (footers):  image2.png -> image1.png
#   BTec logo   (headers):  image1.jpg -> image2.jpg
# The picture's "id"/description/size/position are untouched - only the
# Name (OOXML wp:docPr/@name) changes.

for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            $shp.Name = "image2.jpg"
        }
    }
}

for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            $shp.Name = "image1.png"
        }
    }
}
